# "Add files via upload" — the uploaded workbook has the stray ID value
# ("BABPM2259M") that used to live in cell B18 removed, while the cell
# keeps its existing look (left-aligned / vertically centered, same
# fill & border) — i.e. just the text content goes away.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18").ClearContents()
